$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = -0.06579314222947109
$ws.Range("G2").Value = 0

$ws.Range("B3").Value = -0.07909600764659845
$ws.Range("C3").Value = 0.9807047427925161
$ws.Range("D3").Value = 0.01265768756098197
$ws.Range("E3").Value = 0.007990346254841385
$ws.Range("F3").Value = -0.1018343811709898
$ws.Range("G3").Value = 0.1673002035305011

$ws.Range("B4").Value = -0
$ws.Range("C4").Value = -6.404326871522143
$ws.Range("D4").Value = -0.2425370427083099
$ws.Range("E4").Value = -0.02784059730917719
$ws.Range("F4").Value = 0.9193718120583519
$ws.Range("G4").Value = 0.1298556790611475

$ws.Range("B5").Value = 1.165621356363873
$ws.Range("C5").Value = 5.022181165303445
$ws.Range("D5").Value = 0.4156892904508752
$ws.Range("E5").Value = 0.04685879615052199
$ws.Range("F5").Value = 1.423691301974722
$ws.Range("G5").Value = 0.06494495314546911
